$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
12575,12112,11941,11708,11708,10540,10540,10540,10540,10540,10300,10300,10155,10119,9915,9915,9651,9651,9651,9386,9386,9307,8943,8943,8943,8943,8943,8723,8723,8723,8723,8723,8723,8723,8723,8723,8723,8723,8522,8522,8522,8522,8522,8522,8522,8522,8522,8522,8522,8522,8522,8522,8522,8522,8522,8445,8445,8445,8445,8445,8445,7919,7919,7919,7919,7919,7919,7919,7919,7919,7919,7919,7919,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7295,7295
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $values[$i]
}
